$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.225.20"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = "'1.840.80"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'232.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = "'0.4671"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.07%  '
$ws.Range('D8').Value = "'0.2722"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.89%  '
$ws.Range('D9').Value = "'0.06280"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.61%  '
$ws.Range('D10').Value = "'1.837.40"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('D11').Value = "'0.07420"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = "'16.09"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').Value = "'4.934"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.17%  '
$ws.Range('D14').Value = "'83.72"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.00%  '
$ws.Range('D15').Value = "'0.6199"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.46%  '
$ws.Range('D16').Value = "'30.151.49"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'1.001"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = "'226.05"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.40%  '
$ws.Range('D19').Value = "'0.000007289"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.65%  '
$ws.Range('B20').Value = 'BinanceUSD'
$ws.Range('C20').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D20').Value = "'1.003"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = "'12.33"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.38%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = "'2.078.93"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').Value = "'4.889"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.07%  '
$ws.Range('D24').Value = "'5.853"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.20%  '
$ws.Range('D25').Value = "'9.185"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').Value = "'164.37"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.59%  '
$ws.Range('D27').Value = "'17.75"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.26%  '
$ws.Range('D28').Value = "'1.861"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.38%  '
$ws.Range('D29').Value = "'0.1034"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('D30').Value = "'1.374"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').Value = "'4.073"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.69%  '
$ws.Range('D32').Value = "'3.803"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.72%  '
$ws.Range('D33').Value = "'0.04809"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.54%  '
$ws.Range('D34').Value = "'1.143"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.39%  '
$ws.Range('D35').Value = "'0.7066"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.03%  '
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('D37').Value = "'0.01862"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.01%  '
$ws.Range('D38').Value = "'2.648"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('D39').Value = "'0.8906"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.90%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').Value = "'104.25"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = "'1.913"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.67%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'1.002"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('D43').Value = "'5.514"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('D44').Value = "'0.4010"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.52%  '
$ws.Range('D45').Value = "'7.035"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.73%  '
$ws.Range('D46').Value = "'0.1193"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.88%  '
$ws.Range('D47').Value = "'59.73"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.82%  '
$ws.Range('D48').Value = "'8.591"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.49%  '
$ws.Range('D49').Value = "'32.90"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('D50').Value = "'0.05510"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('E51').Value = '  -4.85%  '

Write-Host "Applied 109 cell updates"
